$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 18.48292350769043
$ws.Range("D2").Value = 190

$ws.Range("C3").Value = 17.84515380859375
$ws.Range("D3").Value = 174

$ws.Range("C4").Value = 17.67802238464355
$ws.Range("D4").Value = 175

$ws.Range("C5").Value = 16.99686050415039
$ws.Range("D5").Value = 193

$ws.Range("C6").Value = 16.57605171203613
$ws.Range("D6").Value = 182
